$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark the "Save File As" button menu bar feature and its section header as Done
$ws.Range("D13").Value = "Yes"
$ws.Range("D14").Value = "Yes"

# Update the selection shown when the workbook is reopened
$ws.Range("A13:A19").Select()
